$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-06-08 Saturday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-06-09 Sunday", 2) | Out-Null

# Update the division problems inside the single table, cell by cell
# (row, col) => new text, addressed positionally to avoid ambiguity
# from values that recur elsewhere in the table.
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "46÷7=6, 4"
$t.Cell(1, 2).Range.Text = "25÷9=2, 7"
$t.Cell(1, 3).Range.Text = "30÷8=3, 6"
$t.Cell(1, 4).Range.Text = "43÷4=10, 3"
$t.Cell(1, 5).Range.Text = "44÷2=22, 0"

$t.Cell(5, 1).Range.Text = "89÷8=11, 1"
$t.Cell(5, 2).Range.Text = "15÷8=1, 7"
$t.Cell(5, 3).Range.Text = "83÷9=9, 2"
$t.Cell(5, 4).Range.Text = "93÷9=10, 3"
$t.Cell(5, 5).Range.Text = "70÷5=14, 0"

$t.Cell(9, 1).Range.Text = "50÷7=7, 1"
$t.Cell(9, 2).Range.Text = "98÷5=19, 3"
$t.Cell(9, 3).Range.Text = "46÷3=15, 1"
$t.Cell(9, 4).Range.Text = "59÷7=8, 3"
$t.Cell(9, 5).Range.Text = "25÷4=6, 1"

$t.Cell(13, 1).Range.Text = "72÷9=8, 0"
$t.Cell(13, 2).Range.Text = "59÷7=8, 3"
$t.Cell(13, 3).Range.Text = "13÷7=1, 6"
$t.Cell(13, 4).Range.Text = "42÷3=14, 0"
$t.Cell(13, 5).Range.Text = "31÷9=3, 4"

$t.Cell(17, 1).Range.Text = "29÷4=7, 1"
$t.Cell(17, 2).Range.Text = "40÷6=6, 4"
$t.Cell(17, 3).Range.Text = "30÷4=7, 2"
$t.Cell(17, 4).Range.Text = "71÷5=14, 1"
$t.Cell(17, 5).Range.Text = "42÷5=8, 2"
